$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.792.54'
$ws.Range("E2").Value = '  +1.48%  '

$ws.Range("D3").Value = '2.804.07'
$ws.Range("E3").Value = '  +1.85%  '

$ws.Range("E4").Value = '  +0.00%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '353.58'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +0.41%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '112.28'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +4.72%  '

$ws.Range("E7").Value = '  +2.38%  '

$ws.Range("E8").Value = '  +0.04%  '

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.624'
$cell.Style = $origStyle
$ws.Range("E9").Value = '  +7.91%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '40.20'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +2.99%  '

$ws.Range("E11").Value = '  -0.41%  '

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0839'
$cell.Style = $origStyle
$ws.Range("E12").Value = '  +1.15%  '

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.90'
$cell.Style = $origStyle
$ws.Range("E13").Value = '  +1.66%  '

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.77'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  +4.24%  '

$ws.Range("D15").Value = '3.245.57'
$ws.Range("E15").Value = '  +2.10%  '

$ws.Range("D16").Value = '2.795.20'
$ws.Range("E16").Value = '  +2.22%  '

$ws.Range("E17").Value = '  +3.16%  '

$ws.Range("D18").Value = '51.806.55'
$ws.Range("E18").Value = '  +1.50%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.61'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  -0.27%  '

$ws.Range("E20").Value = '  +8.41%  '

$ws.Range("E21").Value = '  +4.88%  '

$ws.Range("E22").Value = '  +2.23%  '

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '70.31'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  +1.84%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '267.49'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  +1.65%  '

$ws.Range("E25").Value = '  +1.98%  '

$ws.Range("E26").Value = '  +0.02%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.15'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +1.38%  '

$ws.Range("E28").Value = '  +0.75%  '

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '39.24'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  +13.83%  '

$ws.Range("E30").Value = '  +3.87%  '

$ws.Range("E31").Value = '  +2.26%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '52.29'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  +1.47%  '

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.11'
$cell.Style = $origStyle
$ws.Range("E33").Value = '  +2.19%  '

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0898'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  +9.03%  '

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0449'
$cell.Style = $origStyle
$ws.Range("E35").Value = '  +2.42%  '

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.53'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  +6.43%  '

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  -0.08%  '

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '18.99'
$cell.Style = $origStyle
$ws.Range("E38").Value = '  +3.78%  '

$ws.Range("E39").Value = '  +1.48%  '

$ws.Range("E40").Value = '  +3.71%  '

$ws.Range("E41").Value = '  +2.42%  '

$ws.Range("E42").Value = '  +1.44%  '

$ws.Range("E43").Value = '  +1.61%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '119.71'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  -0.94%  '

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.83'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  -0.70%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.54'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  +10.55%  '

$ws.Range("E47").Value = '  +9.19%  '

$ws.Range("D48").Value = '2.114.50'
$ws.Range("E48").Value = '  +1.35%  '

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.984'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +8.29%  '

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.49'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +1.19%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.37'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +8.28%  '
